$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 43 / 44 swap first (Coin name + Link), then price/volume cells updated below

$ws.Range("D2").Value = '''66.841.39'
$ws.Range("E2").Value = '  +5.39%  '
$ws.Range("D3").Value = '''3.719.28'
$ws.Range("E3").Value = '  +7.42%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''423.92'
$ws.Range("E5").Value = '  +2.12%  '
$ws.Range("D6").Value = '''132.10'
$ws.Range("E6").Value = '  +2.33%  '
$ws.Range("D7").Value = '''3.711.93'
$ws.Range("E7").Value = '  +7.56%  '
$ws.Range("D8").Value = '''0.642'
$ws.Range("E8").Value = '  +2.53%  '
$ws.Range("D9").Value = '''1.00'
$ws.Range("E9").Value = '  +0.05%  '
$ws.Range("D10").Value = '''0.768'
$ws.Range("E10").Value = '  +0.31%  '
$ws.Range("D11").Value = '''0.182'
$ws.Range("E11").Value = '  +13.85%  '
$ws.Range("D12").Value = '''0.0000385'
$ws.Range("E12").Value = '  +63.86%  '
$ws.Range("D13").Value = '''42.89'
$ws.Range("E13").Value = '  +1.97%  '
$ws.Range("E14").Value = '  +5.14%  '
$ws.Range("D15").Value = '''4.307.32'
$ws.Range("E15").Value = '  +7.07%  '
$ws.Range("D17").Value = '''20.92'
$ws.Range("E17").Value = '  +2.97%  '
$ws.Range("D18").Value = '''3.743.88'
$ws.Range("E18").Value = '  +7.68%  '
$ws.Range("D19").Value = '''13.17'
$ws.Range("E19").Value = '  +5.70%  '
$ws.Range("D20").Value = '''1.13'
$ws.Range("E20").Value = '  +3.80%  '
$ws.Range("D21").Value = '''66.857.17'
$ws.Range("E21").Value = '  +5.67%  '
$ws.Range("D22").Value = '''448.21'
$ws.Range("E22").Value = '  -2.51%  '
$ws.Range("D23").Value = '''15.89'
$ws.Range("E23").Value = '  +18.88%  '
$ws.Range("D24").Value = '''90.69'
$ws.Range("E24").Value = '  +0.28%  '
$ws.Range("D25").Value = '''3.19'
$ws.Range("E25").Value = '  -3.26%  '
$ws.Range("D26").Value = '''38.26'
$ws.Range("E26").Value = '  +13.66%  '
$ws.Range("D27").Value = '''10.29'
$ws.Range("E27").Value = '  +1.85%  '
$ws.Range("D28").Value = '''3.31'
$ws.Range("E28").Value = '  +0.54%  '
$ws.Range("D29").Value = '''4.97'
$ws.Range("E29").Value = '  +4.47%  '
$ws.Range("D30").Value = '''12.63'
$ws.Range("E30").Value = '  +2.03%  '
$ws.Range("D31").Value = '''2.79'
$ws.Range("E31").Value = '  +4.54%  '
$ws.Range("E32").Value = '  +7.40%  '
$ws.Range("E33").Value = '  -3.06%  '
$ws.Range("D34").Value = '''42.23'
$ws.Range("E34").Value = '  +5.82%  '
$ws.Range("E35").Value = '  -1.36%  '
$ws.Range("D36").Value = '''0.999'
$ws.Range("E36").Value = '  -0.08%  '
$ws.Range("D37").Value = '''56.42'
$ws.Range("E37").Value = '  -2.43%  '
$ws.Range("D38").Value = '''0.0491'
$ws.Range("E38").Value = '  +0.61%  '
$ws.Range("D39").Value = '''0.0₃0737'
$ws.Range("E39").Value = '  +13.11%  '
$ws.Range("D40").Value = '''3.08'
$ws.Range("E40").Value = '  +32.20%  '
$ws.Range("D41").Value = '''0.147'
$ws.Range("E41").Value = '  +6.74%  '
$ws.Range("D42").Value = '''29.31'
$ws.Range("E42").Value = '  +33.50%  '
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").Value = '''0.996'
$ws.Range("E43").Value = '  -0.30%  '
$ws.Range("B44").Value = 'LidoDAOToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D44").Value = '''3.47'
$ws.Range("E44").Value = '  +4.04%  '
$ws.Range("D45").Value = '''2.14'
$ws.Range("E45").Value = '  +6.95%  '
$ws.Range("D46").Value = '''145.98'
$ws.Range("E46").Value = '  +0.16%  '
$ws.Range("D47").Value = '''2.93'
$ws.Range("E47").Value = '  -4.92%  '
$ws.Range("D48").Value = '''4.40'
$ws.Range("E48").Value = '  -0.84%  '
$ws.Range("D49").Value = '''2.66'
$ws.Range("E49").Value = '  -5.68%  '
$ws.Range("D50").Value = '''0.308'
$ws.Range("E50").Value = '  -3.19%  '
$ws.Range("E51").Value = '  +15.06%  '
